# Update receptor/edge-weight statistics with newly recomputed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 100.4511693333333
$ws.Range("N2").Value = 301.353508
$ws.Range("O2").Value = 0.6815338940941451
$ws.Range("P2").Value = 0.681533894094145
$ws.Range("Q2").Value = 1048.131943941871
$ws.Range("R2").Value = 9433.187495476841
$ws.Range("S2").Value = 0.6618438332212244
$ws.Range("T2").Value = 0.6618438332212244

$ws.Range("N3").Value = 7.755446
$ws.Range("O3").Value = 0.01753953138921768
$ws.Range("P3").Value = 0.01753953138921768
$ws.Range("S3").Value = 0.0170328002585595
$ws.Range("T3").Value = 0.0170328002585595

$ws.Range("M4").Value = 43.839503
$ws.Range("N4").Value = 131.518509
$ws.Range("O4").Value = 0.2974391178622877
$ws.Range("P4").Value = 0.2974391178622877
$ws.Range("Q4").Value = 457.4320419143967
$ws.Range("R4").Value = 4116.888377229569
$ws.Range("S4").Value = 0.2888458631651307
$ws.Range("T4").Value = 0.2888458631651306

$ws.Range("M5").Value = 0.5140156666666666
$ws.Range("N5").Value = 1.542047
$ws.Range("O5").Value = 0.003487456654349595
$ws.Range("P5").Value = 0.003487456654349595
$ws.Range("Q5").Value = 5.363364543145555
$ws.Range("R5").Value = 48.27028088830999
$ws.Range("S5").Value = 0.003386701234243768
$ws.Range("T5").Value = 0.003386701234243768

$ws.Range("M6").Value = 100.4511693333333
$ws.Range("N6").Value = 301.353508
$ws.Range("O6").Value = 0.6815338940941451
$ws.Range("P6").Value = 0.681533894094145
$ws.Range("Q6").Value = 31.182252886792
$ws.Range("R6").Value = 280.640275981128
$ws.Range("S6").Value = 0.01969006087292057
$ws.Range("T6").Value = 0.01969006087292057

$ws.Range("N7").Value = 7.755446
$ws.Range("O7").Value = 0.01753953138921768
$ws.Range("P7").Value = 0.01753953138921768
$ws.Range("R7").Value = 7.222383174636
$ws.Range("S7").Value = 0.000506731130658178
$ws.Range("T7").Value = 0.000506731130658178

$ws.Range("M8").Value = 43.839503
$ws.Range("N8").Value = 131.518509
$ws.Range("O8").Value = 0.2974391178622877
$ws.Range("P8").Value = 0.2974391178622877
$ws.Range("Q8").Value = 13.608746200266
$ws.Range("R8").Value = 122.478715802394
$ws.Range("S8").Value = 0.008593254697157039
$ws.Range("T8").Value = 0.008593254697157039

$ws.Range("M9").Value = 0.5140156666666666
$ws.Range("N9").Value = 1.542047
$ws.Range("O9").Value = 0.003487456654349595
$ws.Range("P9").Value = 0.003487456654349595
$ws.Range("Q9").Value = 0.159561771278
$ws.Range("R9").Value = 1.436055941502
$ws.Range("S9").Value = 0.0001007554201058265
$ws.Range("T9").Value = 0.0001007554201058265
